$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8275113105773926
$ws.Range("B1").Value = 1.265702128410339
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.498253703117371
